# Auto-generated edit script: updates Betfair back/lay odds values
# for rows 2-10 to reflect the latest scraped odds (2025-10-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 2.24
$ws.Cells.Item(2, 7).Value = 2.28
$ws.Cells.Item(2, 8).Value = 3.15
$ws.Cells.Item(2, 9).Value = 3.25
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(2, 12).Value = 1.28
$ws.Cells.Item(2, 14).Value = 5.6
$ws.Cells.Item(2, 15).Value = 1.19
$ws.Cells.Item(2, 16).Value = 2.6
$ws.Cells.Item(2, 17).Value = 1.57
$ws.Cells.Item(2, 18).Value = 1.65
$ws.Cells.Item(2, 19).Value = 2.38
$ws.Cells.Item(2, 21).Value = 2.66
$ws.Cells.Item(2, 22).Value = 1.44
$ws.Cells.Item(2, 23).Value = 1.78
$ws.Cells.Item(2, 24).Value = 24
$ws.Cells.Item(2, 25).Value = 19.5
$ws.Cells.Item(2, 26).Value = 27
$ws.Cells.Item(2, 27).Value = 55
$ws.Cells.Item(2, 30).Value = 14
$ws.Cells.Item(2, 31).Value = 30
$ws.Cells.Item(2, 32).Value = 18
$ws.Cells.Item(2, 34).Value = 14
$ws.Cells.Item(2, 35).Value = 34
$ws.Cells.Item(2, 36).Value = 80
$ws.Cells.Item(2, 37).Value = 20
$ws.Cells.Item(2, 38).Value = 28
$ws.Cells.Item(2, 40).Value = 11.5
# Row 3
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 8).Value = 3.55
$ws.Cells.Item(3, 9).Value = 3.95
$ws.Cells.Item(3, 18).Value = 1.53
$ws.Cells.Item(3, 19).Value = 2.56
$ws.Cells.Item(3, 27).Value = 200
$ws.Cells.Item(3, 30).Value = 19.5
$ws.Cells.Item(3, 39).Value = 300
# Row 4
$ws.Cells.Item(4, 6).Value = 5.6
$ws.Cells.Item(4, 8).Value = 1.68
$ws.Cells.Item(4, 9).Value = 1.71
$ws.Cells.Item(4, 11).Value = 4.6
$ws.Cells.Item(4, 12).Value = 1.37
$ws.Cells.Item(4, 13).Value = 1.05
$ws.Cells.Item(4, 14).Value = 4.1
$ws.Cells.Item(4, 17).Value = 1.81
$ws.Cells.Item(4, 18).Value = 1.39
$ws.Cells.Item(4, 19).Value = 3.1
$ws.Cells.Item(4, 21).Value = 2.08
$ws.Cells.Item(4, 22).Value = 2.4
$ws.Cells.Item(4, 23).Value = 1.2
$ws.Cells.Item(4, 26).Value = 10.5
$ws.Cells.Item(4, 27).Value = 17
$ws.Cells.Item(4, 29).Value = 10.5
$ws.Cells.Item(4, 30).Value = 10.5
$ws.Cells.Item(4, 31).Value = 18
$ws.Cells.Item(4, 35).Value = 980
$ws.Cells.Item(4, 37).Value = 380
$ws.Cells.Item(4, 39).Value = 580
$ws.Cells.Item(4, 41).Value = 10
# Row 5
$ws.Cells.Item(5, 6).Value = 1.88
$ws.Cells.Item(5, 7).Value = 2.08
$ws.Cells.Item(5, 8).Value = 4.3
$ws.Cells.Item(5, 9).Value = 5.2
$ws.Cells.Item(5, 10).Value = 3.35
$ws.Cells.Item(5, 11).Value = 3.85
$ws.Cells.Item(5, 14).Value = 3.1
$ws.Cells.Item(5, 16).Value = 1.64
$ws.Cells.Item(5, 17).Value = 2.02
$ws.Cells.Item(5, 18).Value = 1.22
$ws.Cells.Item(5, 19).Value = 3.55
$ws.Cells.Item(5, 20).Value = 1.64
$ws.Cells.Item(5, 21).Value = 1.62
$ws.Cells.Item(5, 22).Value = 1.24
$ws.Cells.Item(5, 23).Value = 1.92
$ws.Cells.Item(5, 24).Value = 17
$ws.Cells.Item(5, 28).Value = 10.5
$ws.Cells.Item(5, 29).Value = 11
# Row 6
$ws.Cells.Item(6, 6).Value = 3.05
$ws.Cells.Item(6, 7).Value = 3.15
$ws.Cells.Item(6, 8).Value = 2.76
$ws.Cells.Item(6, 9).Value = 2.82
$ws.Cells.Item(6, 11).Value = 3.15
$ws.Cells.Item(6, 14).Value = 2.62
$ws.Cells.Item(6, 15).Value = 1.58
$ws.Cells.Item(6, 16).Value = 1.54
$ws.Cells.Item(6, 17).Value = 2.72
$ws.Cells.Item(6, 21).Value = 1.78
$ws.Cells.Item(6, 23).Value = 1.47
$ws.Cells.Item(6, 24).Value = 8.6
$ws.Cells.Item(6, 25).Value = 8.199999999999999
$ws.Cells.Item(6, 26).Value = 17
$ws.Cells.Item(6, 27).Value = 50
$ws.Cells.Item(6, 28).Value = 8.800000000000001
$ws.Cells.Item(6, 29).Value = 7.2
$ws.Cells.Item(6, 30).Value = 14.5
$ws.Cells.Item(6, 31).Value = 46
$ws.Cells.Item(6, 32).Value = 19
$ws.Cells.Item(6, 34).Value = 26
$ws.Cells.Item(6, 35).Value = 170
$ws.Cells.Item(6, 36).Value = 150
$ws.Cells.Item(6, 37).Value = 50
$ws.Cells.Item(6, 38).Value = 210
$ws.Cells.Item(6, 39).Value = 240
$ws.Cells.Item(6, 40).Value = 210
$ws.Cells.Item(6, 41).Value = 55
# Row 7
$ws.Cells.Item(7, 7).Value = 1.68
$ws.Cells.Item(7, 8).Value = 5.7
$ws.Cells.Item(7, 10).Value = 3.75
$ws.Cells.Item(7, 14).Value = 1.94
$ws.Cells.Item(7, 15).Value = 1.24
$ws.Cells.Item(7, 16).Value = 1.94
$ws.Cells.Item(7, 17).Value = 1.56
$ws.Cells.Item(7, 19).Value = 1.56
$ws.Cells.Item(7, 23).Value = 2.46
# Row 8
$ws.Cells.Item(8, 6).Value = 1.44
$ws.Cells.Item(8, 17).Value = 2.02
$ws.Cells.Item(8, 19).Value = 3.65
$ws.Cells.Item(8, 25).Value = 28
$ws.Cells.Item(8, 27).Value = 540
$ws.Cells.Item(8, 31).Value = 260
$ws.Cells.Item(8, 41).Value = 450
# Row 9
$ws.Cells.Item(9, 6).Value = 3
$ws.Cells.Item(9, 7).Value = 3.1
$ws.Cells.Item(9, 8).Value = 2.78
$ws.Cells.Item(9, 9).Value = 2.88
$ws.Cells.Item(9, 10).Value = 3.1
$ws.Cells.Item(9, 11).Value = 3.15
$ws.Cells.Item(9, 14).Value = 2.72
$ws.Cells.Item(9, 16).Value = 1.55
$ws.Cells.Item(9, 17).Value = 2.66
$ws.Cells.Item(9, 21).Value = 1.79
$ws.Cells.Item(9, 22).Value = 1.53
$ws.Cells.Item(9, 23).Value = 1.48
# Row 10
$ws.Cells.Item(10, 6).Value = 1.38
$ws.Cells.Item(10, 7).Value = 1.4
$ws.Cells.Item(10, 8).Value = 10.5
$ws.Cells.Item(10, 9).Value = 14
$ws.Cells.Item(10, 10).Value = 5.1
$ws.Cells.Item(10, 14).Value = 3.4
$ws.Cells.Item(10, 15).Value = 1.38
$ws.Cells.Item(10, 16).Value = 1.84
$ws.Cells.Item(10, 17).Value = 2.08
$ws.Cells.Item(10, 18).Value = 1.31
$ws.Cells.Item(10, 19).Value = 3.9
$ws.Cells.Item(10, 20).Value = 2.5
$ws.Cells.Item(10, 21).Value = 1.58
$ws.Cells.Item(10, 23).Value = 3.45
$ws.Cells.Item(10, 25).Value = 29
$ws.Cells.Item(10, 27).Value = 710
$ws.Cells.Item(10, 30).Value = 50
$ws.Cells.Item(10, 31).Value = 330
$ws.Cells.Item(10, 34).Value = 42
$ws.Cells.Item(10, 35).Value = 270
$ws.Cells.Item(10, 39).Value = 340
$ws.Cells.Item(10, 41).Value = 600
